# Apply the BrazilResponsiveness.xlsx update:
#  1. Extend the Brazil responsiveness rows (1,3,5) from column V out to column AE
#     (i.e. from a 20-yr run, 2001-2020, to a 30-yr run, 2001-2030), repeating the
#     last value across the new columns.
#  2. Add new input rows for US/EU27 yield-variation modifiers (rows 7,9,11,13),
#     each defaulted to 0 across the full B:AE run.
#  3. Add new rows controlling the China (CHIHKG) western-diet transition
#     mechanism: on/off switch (row 15), start year (row 17), and the
#     exponential-averaging transition time (row 19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend existing rows 1, 3, 5 from column V to column AE ------------
# Copy column V's formatting (keeps the orange/"s=1" style on row 1) then
# stamp the same value across the newly-added columns, matching the way the
# existing B:V runs already look (constant value repeated to the right).

$ws.Range("V1").Copy() | Out-Null
$ws.Range("W1:AE1").PasteSpecial(-4122) | Out-Null
$ws.Range("W1:AE1").Value = 1.75

$ws.Range("V3").Copy() | Out-Null
$ws.Range("W3:AE3").PasteSpecial(-4122) | Out-Null
$ws.Range("W3:AE3").Value = 0.609

$ws.Range("V5").Copy() | Out-Null
$ws.Range("W5:AE5").PasteSpecial(-4122) | Out-Null
$ws.Range("W5:AE5").Value = 0.01

# --- 2. New yield-variation-modifier rows (USA / EU27, Maize / OilCrop) ----

$ws.Range("A7").Value = "USA.C_yield_variation_PCT[Maize]"
$ws.Range("B7:AE7").Value = 0

$ws.Range("A9").Value = "USA.C_yield_variation_PCT[OilCrop]"
$ws.Range("B9:AE9").Value = 0

$ws.Range("A11").Value = "EU27.C_yield_variation_PCT[Maize]"
$ws.Range("B11:AE11").Value = 0

$ws.Range("A13").Value = "EU27.C_yield_variation_PCT[OilCrop]"
$ws.Range("B13:AE13").Value = 0

# --- 3. New China (CHIHKG) western-diet transition controls ---------------

$ws.Range("A15").Value = "CHIHKG.TransitionToWesternDiet"
$ws.Range("B15").Value = 0

$ws.Range("A17").Value = "CHIHKG.DietTransition Start"
$ws.Range("B17").Value = 2020

$ws.Range("A19").Value = "CHIHKG.DietTransitionTime"
$ws.Range("B19").Value = 5

# --- Selection matches where the author's cursor ended up on save ---------
$ws.Range("B16").Select() | Out-Null
